$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 3423.5789
$ws.Range("J32").Value = 3276.5
$ws.Range("L32").Value = 3276.5
$ws.Range("N32").Value = -3928.5
$ws.Range("H33").Value = 556.1667
$ws.Range("I33").Value = 510.5
$ws.Range("K33").Value = 510.5
$ws.Range("M33").Value = -281.5
$ws.Range("H39").Value = 1958.8182
$ws.Range("J39").Value = 10249.5
$ws.Range("L39").Value = 30748.5
$ws.Range("N39").Value = -31340.5
$ws.Range("H96").Value = 489.25
$ws.Range("I96").Value = 511
$ws.Range("K96").Value = 1533
$ws.Range("M96").Value = -160
$ws.Range("H113").Value = 4618.7334
$ws.Range("I113").Value = 2902.2144
$ws.Range("K113").Value = 2902.2144
$ws.Range("M113").Value = 351.7856000000002
$ws.Range("H129").Value = 1605.4546
$ws.Range("I129").Value = 1358.5
$ws.Range("J129").Value = 2264
$ws.Range("K129").Value = 4075.5
$ws.Range("L129").Value = 6792
$ws.Range("M129").Value = 924.5
$ws.Range("N129").Value = -16792
$ws.Range("H138").Value = 3300.7083
$ws.Range("I138").Value = 3142.1667
$ws.Range("J138").Value = 3459.25
$ws.Range("K138").Value = 9426.500100000001
$ws.Range("L138").Value = 10377.75
$ws.Range("M138").Value = -4286.500100000001
$ws.Range("N138").Value = -20657.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2013.4865
$ws.Range("I32").Value = 1019.9355
$ws.Range("K32").Value = 1019.9355
$ws.Range("M32").Value = -732.9355
$ws.Range("H74").Value = 1888.08
$ws.Range("I74").Value = 1331.8889
$ws.Range("J74").Value = 3318.2856
$ws.Range("K74").Value = 1331.8889
$ws.Range("L74").Value = 3318.2856
$ws.Range("M74").Value = -457.8888999999999
$ws.Range("N74").Value = -5066.2856
$ws.Range("H77").Value = 1888.08
$ws.Range("I77").Value = 1331.8889
$ws.Range("J77").Value = 3318.2856
$ws.Range("K77").Value = 6659.4445
$ws.Range("L77").Value = 16591.428
$ws.Range("M77").Value = -2291.4445
$ws.Range("N77").Value = -25327.428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3782.8572
$ws.Range("I58").Value = 2582.625
$ws.Range("K58").Value = 2582.625
$ws.Range("M58").Value = -2379.625
$ws.Range("H68").Value = 39997
$ws.Range("J68").Value = 39997
$ws.Range("L68").Value = 39997
$ws.Range("N68").Value = -41495
$ws.Range("H71").Value = 39997
$ws.Range("J71").Value = 39997
$ws.Range("L71").Value = 119991
$ws.Range("N71").Value = -127479
$ws.Range("H92").Value = 25112.25
$ws.Range("J92").Value = 25112.25
$ws.Range("L92").Value = 25112.25
$ws.Range("N92").Value = -30104.25
$ws.Range("H99").Value = 2103.4
$ws.Range("I99").Value = 1955.2858
$ws.Range("K99").Value = 1955.2858
$ws.Range("M99").Value = -457.2858000000001
$ws.Range("H126").Value = 2103.4
$ws.Range("I126").Value = 1955.2858
$ws.Range("K126").Value = 5865.857400000001
$ws.Range("M126").Value = -3395.857400000001
$ws.Range("H132").Value = 3221.238
$ws.Range("I132").Value = 2822.5625
$ws.Range("J132").Value = 4497
$ws.Range("K132").Value = 8467.6875
$ws.Range("L132").Value = 13491
$ws.Range("M132").Value = -5937.6875
$ws.Range("N132").Value = -18551
$ws.Range("H134").Value = 5320.782
$ws.Range("I134").Value = 4758.113
$ws.Range("J134").Value = 7501.125
$ws.Range("K134").Value = 14274.339
$ws.Range("L134").Value = 22503.375
$ws.Range("M134").Value = -11739.339
$ws.Range("N134").Value = -27573.375
$ws.Range("H136").Value = 3782.8572
$ws.Range("I136").Value = 2582.625
$ws.Range("K136").Value = 7747.875
$ws.Range("M136").Value = -5197.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 284.76923
$ws.Range("I23").Value = 199.57143
$ws.Range("J23").Value = 384.16666
$ws.Range("K23").Value = 598.71429
$ws.Range("L23").Value = 1152.49998
$ws.Range("M23").Value = -363.71429
$ws.Range("N23").Value = -1622.49998
$ws.Range("H112").Value = 16571
$ws.Range("I112").Value = 7998.5
$ws.Range("J112").Value = 20000
$ws.Range("K112").Value = 23995.5
$ws.Range("L112").Value = 60000
$ws.Range("M112").Value = -22887.5
$ws.Range("N112").Value = -62216
$ws.Range("H115").Value = 7999.5
$ws.Range("I115").Value = 1999
$ws.Range("J115").Value = 14000
$ws.Range("K115").Value = 5997
$ws.Range("L115").Value = 42000
$ws.Range("M115").Value = -4822
$ws.Range("N115").Value = -44350
$ws.Range("H132").Value = 2057
$ws.Range("I132").Value = 1057.6
$ws.Range("K132").Value = 9518.4
$ws.Range("M132").Value = -6988.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1584.3334
$ws.Range("I122").Value = 1123.6666
$ws.Range("K122").Value = 3370.9998
$ws.Range("M122").Value = -920.9998000000001
$ws.Range("H132").Value = 1863.5217
$ws.Range("I132").Value = 1748.5
$ws.Range("K132").Value = 5245.5
$ws.Range("M132").Value = -2715.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5333.3335
$ws.Range("J22").Value = 5333.3335
$ws.Range("L22").Value = 5333.3335
$ws.Range("N22").Value = -5923.3335
$ws.Range("H27").Value = 5333.3335
$ws.Range("J27").Value = 5333.3335
$ws.Range("L27").Value = 5333.3335
$ws.Range("N27").Value = -5547.3335
$ws.Range("H132").Value = 2617.3428
$ws.Range("I132").Value = 1659.0588
$ws.Range("K132").Value = 4977.1764
$ws.Range("M132").Value = -2447.1764
$ws.Range("H136").Value = 2212.9697
$ws.Range("I136").Value = 1330.2222
$ws.Range("K136").Value = 3990.6666
$ws.Range("M136").Value = -1440.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 5500
$ws.Range("I34").Value = 5500
$ws.Range("K34").Value = 5500
$ws.Range("M34").Value = -5297
$ws.Range("H95").Value = 30000
$ws.Range("J95").Value = 30000
$ws.Range("L95").Value = 30000
$ws.Range("N95").Value = -35492
$ws.Range("H96").Value = 3121.1333
$ws.Range("I96").Value = 2498.75
$ws.Range("J96").Value = 3347.4546
$ws.Range("K96").Value = 2498.75
$ws.Range("L96").Value = 3347.4546
$ws.Range("M96").Value = -1125.75
$ws.Range("N96").Value = -6093.4546
$ws.Range("H97").Value = 29973.8
$ws.Range("J97").Value = 29973.8
$ws.Range("L97").Value = 29973.8
$ws.Range("N97").Value = -31955.8
$ws.Range("H132").Value = 2527.3215
$ws.Range("I132").Value = 2448.7917
$ws.Range("K132").Value = 7346.375100000001
$ws.Range("M132").Value = -4816.375100000001
